$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.095.59"
$ws.Range("E2").Value = "  -2.06%  "

$ws.Range("D3").Value = "3.118.17"
$ws.Range("E3").Value = "  -0.69%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.36%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").Value = "3.109.78"
$ws.Range("E8").Value = "  -0.84%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.515"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.82%  "

$ws.Range("E10").Value = "  -3.72%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.28"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.04%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.455"
$ws.Range("D12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000246"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.56%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.80%  "

$ws.Range("D15").Value = "3.629.99"
$ws.Range("E15").Value = "  -0.68%  "

$ws.Range("E16").Value = "  +2.39%  "

$ws.Range("D17").Value = "62.974.56"
$ws.Range("E17").Value = "  -2.19%  "

$ws.Range("D18").Value = "3.119.83"
$ws.Range("E18").Value = "  -0.35%  "

$ws.Range("E19").Value = "  -2.89%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "472.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.07%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.45%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.696"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.80%  "

$ws.Range("E23").Value = "  -0.25%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "86.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.67%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.93%  "

$ws.Range("E26").Value = "  -0.01%  "

$ws.Range("E27").Value = "  -1.53%  "

$ws.Range("B28").Value = "NEARProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.95"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.74%  "

$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.91"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.44%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.29%  "

$ws.Range("E31").Value = "  +0.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.58%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.108"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.94%  "

$ws.Range("E34").Value = "  -4.16%  "

$ws.Range("E35").Value = "  -2.54%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.79"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.10%  "

$ws.Range("E37").Value = "  -0.68%  "

$ws.Range("D38").Value = "0.0₃0703"
$ws.Range("E38").Value = "  -8.17%  "

$ws.Range("E39").Value = "  -1.85%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "419.70"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.83%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.19"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.72%  "

$ws.Range("D42").Value = "2.899.50"
$ws.Range("E42").Value = "  +0.81%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -11.47%  "

$ws.Range("E44").Value = "  -6.27%  "

$ws.Range("E45").Value = "  +1.35%  "

$ws.Range("E46").Value = "  +0.01%  "

$ws.Range("E47").Value = "  -5.46%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.45"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.87%  "

$ws.Range("E49").Value = "  -0.40%  "

$ws.Range("E50").Value = "  -5.66%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "119.56"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.01%  "

Write-Host "Updated cryptos list"
